# Trade #5 closed at 2026-02-16 22:52:26 - base_strategy DOWN +0.000%
#
# Appends a new trade row (row 6) to both the "All Trades" sheet and the
# "base_strategy" sheet, mirroring the existing OPEN trade rows above it.
#
# We copy the last existing data row (row 5) down into row 6 first -
# this is important because plainly assigning a date-look-alike string
# like "2026-02-16" to a cell's .Value causes Excel's COM layer to parse
# it into a date serial number. Row 5's Date cell is already stored as
# plain text, so copying it preserves the text representation. We then
# only overwrite the two cells (Trade # and Time) that actually differ
# for the new trade.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Clone row 5 (the previous trade row) into the new row 6.
    $ws.Range("A5:Q5").Copy($ws.Range("A6:Q6"))

    # Update the fields that differ for this new trade.
    $ws.Cells.Item(6, 1).Value = 5              # Trade #
    $ws.Cells.Item(6, 3).Value = "22:52:26"     # Time
}
